$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the two comments (with author + rich formatting) before the row
# shift moves their anchors.
$commentShader = $ws.Range("B10").Comment.Text()
$commentShaderAuthor = $ws.Range("B10").Comment.Author
$commentError = $ws.Range("B16").Comment.Text()
$commentErrorAuthor = $ws.Range("B16").Comment.Author

# The "Point lights - specular lighting" task (row 2) was removed; deleting
# the whole row shifts every row below it up by one and drops the now-unused
# shared string automatically.
$ws.Rows(2).Delete()

# Two of the shifted rows also received independent value edits beyond the
# plain shift (Estimate column on the new row 2 and row 3).
$ws.Range("C2").Value = 4
$ws.Range("C3").Value = 2

# Re-create the comments on their new (shifted-up) anchor cells.
$ws.Range("B10").Comment.Delete()
$ws.Range("B16").Comment.Delete()

$newShaderComment = $ws.Range("B9").AddComment($commentShader)
$newShaderComment.Author = $commentShaderAuthor

$newErrorComment = $ws.Range("B15").AddComment($commentError)
$newErrorComment.Author = $commentErrorAuthor

# Match the author's final selection.
$ws.Range("C3").Select() | Out-Null
